# Title slide: "Graph database " -> "Vince Graph database "
# The original single run ("Graph database ") is split into two runs:
#   run 1: "Vince Graph "
#   run 2: "database "
# Achieved by replacing the leading "Graph " substring with "Vince Graph ".

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

$titleShape = $null
for ($i = 1; $i -le $s.Shapes.Count; $i++) {
    $sh = $s.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.TextFrame.TextRange.Text -like "Graph database*") {
            $titleShape = $sh
            break
        }
    }
}

if ($titleShape -eq $null) {
    $titleShape = $s.Shapes.Item(1)
}

$tr = $titleShape.TextFrame.TextRange
$lead = $tr.Characters(1, 6)   # "Graph "
$lead.Text = "Vince Graph "
